$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Update status text "In Translation" -> "Ready for handoff"
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# Update timestamps
$overview.Range("G2").Value = "2016-08-31 07:07:33"
$zhcn.Range("H2").Value = "2016-08-31 07:07:27"
$dede.Range("H2").Value = "2016-08-31 07:07:33"

# Autofit the status columns since the new text is wider
$overview.Range("E:F").EntireColumn.AutoFit()
$zhcn.Range("C:C").EntireColumn.AutoFit()
$dede.Range("C:C").EntireColumn.AutoFit()
